$d = $word.ActiveDocument

# Append a "DICAS DE PREPARO" section (with two tips, separated by blank
# paragraphs) right after the "esfriar." paragraph that closes the
# "MODO DE PREPARO" section, and before the document's trailing blank
# paragraph.
#
# Use Find/Replace on the exact "esfriar." paragraph text: it is matched
# once, and the replacement re-uses it as a prefix followed by `r`-joined
# new paragraphs, so every inserted chunk of text lands in its own new
# paragraph without disturbing the existing "esfriar." paragraph or the
# document's final (pre-existing) empty paragraph.

$oldText = "esfriar."
$newText = "esfriar.`r" + `
    "DICAS DE PREPARO:`r" + `
    "`r" + `
    "- Para um bolo mais fofo, peneire a farinha de trigo.`r" + `
    "`r" + `
    "- Você pode adicionar nozes picadas à massa para um toque especial."

$found = $d.Content.Find.Execute(
    $oldText, $false, $true, $false, $false, $false,
    $true, 1, $false, $newText, 2)

Write-Output $found
